$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 179.1580256666667
$ws.Range("H2").Value = 537.4740770000001
$ws.Range("I2").Value = 0.3468013736386751
$ws.Range("J2").Value = 0.3468013736386751
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.11318066666666
$ws.Range("N2").Value = 66.33954199999999
$ws.Range("O2").Value = 0.1993293533530854
$ws.Range("P2").Value = 0.1993293533530854
$ws.Range("Q2").Value = 3961.753789450304
$ws.Range("R2").Value = 35655.78410505273
$ws.Range("S2").Value = 0.06912769354935887
$ws.Range("T2").Value = 0.06912769354935885
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 179.1580256666667
$ws.Range("H3").Value = 537.4740770000001
$ws.Range("I3").Value = 0.3468013736386751
$ws.Range("J3").Value = 0.3468013736386751
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.622575
$ws.Range("N3").Value = 22.867725
$ws.Range("O3").Value = 0.06871028498970018
$ws.Range("P3").Value = 0.06871028498970018
$ws.Range("Q3").Value = 1365.645487496092
$ws.Range("R3").Value = 12290.80938746483
$ws.Range("S3").Value = 0.02382882121753286
$ws.Range("T3").Value = 0.02382882121753286
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 179.1580256666667
$ws.Range("H4").Value = 537.4740770000001
$ws.Range("I4").Value = 0.3468013736386751
$ws.Range("J4").Value = 0.3468013736386751
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 81.20214833333334
$ws.Range("N4").Value = 243.606445
$ws.Range("O4").Value = 0.7319603616572145
$ws.Range("P4").Value = 0.7319603616572145
$ws.Range("Q4").Value = 14548.01657529181
$ws.Range("R4").Value = 130932.1491776263
$ws.Range("S4").Value = 0.2538448588717834
$ws.Range("T4").Value = 0.2538448588717834
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 239.807332
$ws.Range("H5").Value = 719.421996
$ws.Range("I5").Value = 0.4642019905988459
$ws.Range("J5").Value = 0.4642019905988459
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.11318066666666
$ws.Range("N5").Value = 66.33954199999999
$ws.Range("O5").Value = 0.1993293533530854
$ws.Range("P5").Value = 0.1993293533530854
$ws.Range("Q5").Value = 5302.902857707314
$ws.Range("R5").Value = 47726.12571936583
$ws.Range("S5").Value = 0.09252908261128298
$ws.Range("T5").Value = 0.09252908261128298
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 239.807332
$ws.Range("H6").Value = 719.421996
$ws.Range("I6").Value = 0.4642019905988459
$ws.Range("J6").Value = 0.4642019905988459
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.622575
$ws.Range("N6").Value = 22.867725
$ws.Range("O6").Value = 0.06871028498970018
$ws.Range("P6").Value = 0.06871028498970018
$ws.Range("Q6").Value = 1827.9493737199
$ws.Range("R6").Value = 16451.5443634791
$ws.Range("S6").Value = 0.03189545106683283
$ws.Range("T6").Value = 0.03189545106683283
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 239.807332
$ws.Range("H7").Value = 719.421996
$ws.Range("I7").Value = 0.4642019905988459
$ws.Range("J7").Value = 0.4642019905988459
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 81.20214833333334
$ws.Range("N7").Value = 243.606445
$ws.Range("O7").Value = 0.7319603616572145
$ws.Range("P7").Value = 0.7319603616572145
$ws.Range("Q7").Value = 19472.87054448491
$ws.Range("R7").Value = 175255.8349003642
$ws.Range("S7").Value = 0.3397774569207301
$ws.Range("T7").Value = 0.3397774569207301
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 97.63589966666666
$ws.Range("H8").Value = 292.907699
$ws.Range("I8").Value = 0.1889966357624789
$ws.Range("J8").Value = 0.1889966357624789
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 22.11318066666666
$ws.Range("N8").Value = 66.33954199999999
$ws.Range("O8").Value = 0.1993293533530854
$ws.Range("P8").Value = 0.1993293533530854
$ws.Range("Q8").Value = 2159.040288881539
$ws.Range("R8").Value = 19431.36259993385
$ws.Range("S8").Value = 0.03767257719244354
$ws.Range("T8").Value = 0.03767257719244354
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 97.63589966666666
$ws.Range("H9").Value = 292.907699
$ws.Range("I9").Value = 0.1889966357624789
$ws.Range("J9").Value = 0.1889966357624789
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.622575
$ws.Range("N9").Value = 22.867725
$ws.Range("O9").Value = 0.06871028498970018
$ws.Range("P9").Value = 0.06871028498970018
$ws.Range("Q9").Value = 744.2369679016416
$ws.Range("R9").Value = 6698.132711114775
$ws.Range("S9").Value = 0.01298601270533449
$ws.Range("T9").Value = 0.01298601270533449
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 97.63589966666666
$ws.Range("H10").Value = 292.907699
$ws.Range("I10").Value = 0.1889966357624789
$ws.Range("J10").Value = 0.1889966357624789
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 81.20214833333334
$ws.Range("N10").Value = 243.606445
$ws.Range("O10").Value = 0.7319603616572145
$ws.Range("P10").Value = 0.7319603616572145
$ws.Range("Q10").Value = 7928.244807391117
$ws.Range("R10").Value = 71354.20326652005
$ws.Range("S10").Value = 0.1383380458647009
$ws.Range("T10").Value = 0.1383380458647009
